# Apply the updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as literal text in the sheet (e.g. "62.979.24").
# Some updated prices are plain decimals (e.g. "146.32") which Excel would
# otherwise auto-convert to a number on assignment, so force those cells to
# Text format first to keep them as literal strings, matching the source data.

$ws.Range("D2").Value = "62.979.24"
$ws.Range("E2").Value = "  +3.12%  "
$ws.Range("D3").Value = "2.453.79"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.12"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.32"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "2.453.37"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("E14").Value = "  +7.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +5.55%  "
$ws.Range("D16").Value = "2.896.41"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").Value = "62.866.04"
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").Value = "2.447.50"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.08"
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.49"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.08"
$ws.Range("E23").Value = "  +8.53%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.40"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "650.50"
$ws.Range("E26").Value = "  +10.65%  "
$ws.Range("E27").Value = "  +17.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.57"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.22"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("E32").Value = "  +6.55%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("B34").Value = "BabyDogeCoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D34").Value = "0.0₆0435"
$ws.Range("E34").Value = "  +54.96%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.78"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  +6.20%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.98"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").Value = "  +10.55%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.54"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.00"
$ws.Range("E47").Value = "  +27.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.03"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.64"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.74"
$ws.Range("E50").Value = "  +5.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.607"
$ws.Range("E51").Value = "  +2.53%  "
